# Update countries & provincias Spain
# Refresh the "Pais" sheet with the next data pull (17 Oct 23:54 -> 18 Oct 01:11):
#   - bump the "last updated" timestamp banner
#   - update case counters for the countries whose totals changed
#   - three countries (Chequia, Nigeria, Angola) overtook their neighbour
#     in the ranking, so those row pairs swap country name + stats

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Banner: "Datos actualizados a ..." ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 01:11"

# --- Plain stat refreshes (country/rank unchanged) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 8340243
$ws.Range("C4").Value = 51810
$ws.Range("D4").Value = 5427299
$ws.Range("E4").Value = 2688680
$ws.Range("G4").Value = 620
$ws.Range("H4").Value = 224264

# Row 9: Argentina
$ws.Range("B9").Value = 979119
$ws.Range("C9").Value = 13510
$ws.Range("D9").Value = 791174
$ws.Range("E9").Value = 161838
$ws.Range("G9").Value = 384
$ws.Range("H9").Value = 26107

# Row 10: Colombia
$ws.Range("B10").Value = 952371
$ws.Range("C10").Value = 7017
$ws.Range("D10").Value = 847467
$ws.Range("E10").Value = 76101
$ws.Range("G10").Value = 187
$ws.Range("H10").Value = 28803

# Row 12: Peru
$ws.Range("B12").Value = 865549
$ws.Range("C12").Value = 3132
$ws.Range("D12").Value = 774356
$ws.Range("E12").Value = 57491
$ws.Range("G12").Value = 54
$ws.Range("H12").Value = 33702

# Row 40: Panama
$ws.Range("B40").Value = 124107
$ws.Range("C40").Value = 609
$ws.Range("D40").Value = 100348
$ws.Range("E40").Value = 21202
$ws.Range("G40").Value = 11
$ws.Range("H40").Value = 2557

# Row 46: Egipto
$ws.Range("B46").Value = 105297
$ws.Range("C46").Value = 138
$ws.Range("D46").Value = 98157
$ws.Range("E46").Value = 1031
$ws.Range("G46").Value = 10
$ws.Range("H46").Value = 6109

# Row 58: Suiza (only Casos activos / Recuperados moved)
$ws.Range("D58").Value = 50600
$ws.Range("E58").Value = 21700

# Row 84: Bulgaria
$ws.Range("B84").Value = 29108
$ws.Range("C84").Value = 603
$ws.Range("D84").Value = 16912
$ws.Range("E84").Value = 11228
$ws.Range("G84").Value = 10
$ws.Range("H84").Value = 968

# Row 118: Mauritania (Casos criticos / Muertes hoy / Muertes untouched)
$ws.Range("B118").Value = 7607
$ws.Range("C118").Value = 4
$ws.Range("D118").Value = 7342
$ws.Range("E118").Value = 102

# Row 134: Guinea Ecuatorial (Casos activos untouched)
$ws.Range("B134").Value = 5070
$ws.Range("C134").Value = 2
$ws.Range("E134").Value = 33

# --- Ranking swaps: country moves up, its neighbour slides down one row ---

# Rows 34/35: Chequia overtakes Polonia
$ws.Range("A34").Value = "Chequia"
$ws.Range("B34").Value = 168827
$ws.Range("C34").Value = 8715
$ws.Range("D34").Value = 68945
$ws.Range("E34").Value = 98530
$ws.Range("G34").Value = 69
$ws.Range("H34").Value = 1352

$ws.Range("A35").Value = "Polonia"
$ws.Range("B35").Value = 167230
$ws.Range("C35").Value = 9622
$ws.Range("D35").Value = 90162
$ws.Range("E35").Value = 73544
$ws.Range("G35").Value = 84
$ws.Range("H35").Value = 3524

# Rows 63/64: Nigeria overtakes Libano
$ws.Range("A63").Value = "Nigeria"
$ws.Range("B63").Value = 61307
$ws.Range("C63").Value = 113
$ws.Range("D63").Value = 56557
$ws.Range("E63").Value = 3627
$ws.Range("G63").Value = 4
$ws.Range("H63").Value = 1123

$ws.Range("A64").Value = "Libano"
$ws.Range("B64").Value = 61284
$ws.Range("C64").Value = 1171
$ws.Range("D64").Value = 27197
$ws.Range("E64").Value = 33570
$ws.Range("G64").Value = 8
$ws.Range("H64").Value = 517

# Rows 119/120: Angola overtakes Lituania
$ws.Range("A119").Value = "Angola"
$ws.Range("B119").Value = 7462
$ws.Range("C119").Value = 240
$ws.Range("D119").Value = 3022
$ws.Range("E119").Value = 4199
$ws.Range("G119").Value = 7
$ws.Range("H119").Value = 241

$ws.Range("A120").Value = "Lituania"
$ws.Range("B120").Value = 7269
$ws.Range("C120").Value = 228
$ws.Range("D120").Value = 3097
$ws.Range("E120").Value = 4059
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 113
